$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data pull.
# D-column prices are plain numeric-looking text (e.g. "20.60", "62.689.01") in the
# source sheet, so force text formatting before assignment to avoid Excel silently
# coercing them into numbers (which would drop meaningful trailing zeros or choke on
# multi-dot thousand-separated values).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.689.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.469.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("E7").Value = "  -0.61%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("E10").Value = "  +9.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.72%  "

$ws.Range("E13").Value = "  +3.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.023.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.505.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.667.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "463.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +16.38%  "

$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("E28").Value = "  +0.13%  "

$ws.Range("E29").Value = "  -1.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("E31").Value = "  -3.22%  "

$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.25%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.18%  "

$ws.Range("E37").Value = "  -2.81%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  +3.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("E42").Value = "  +5.73%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("E45").Value = "  +3.65%  "

$ws.Range("E46").Value = "  +3.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0571"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +34.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.140"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.43%  "
